$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Result / Actions" text for the first Activity row, plus a new
#     second row of content (same Activity, a different result) ---
$ws.Range("C2").Value = "Implement greedy algorithm, the bag problem of the dynamic plan;"
$ws.Range("C3").Value = "Go through the concept of several classifier, like KNN, Naïve Bayes, and introduction of other common and important algorithms, like MapReduce (distributed algorithm)"

# Row 3 used to be its own "No 2" entry; now it is just the second result
# line for the same Activity, so clear the old No/Activity values before
# merging the No (A) and Activity (B) cells across rows 2-3.
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()

# Make column B match column A's centered style so the merged cell looks
# consistent (it was left-aligned before).
$ws.Range("B2:B3").HorizontalAlignment = -4108

$ws.Range("B2:B3").Merge()
$ws.Range("A2:A3").Merge()

# Keep both halves of the now taller activity block the same height.
$ws.Rows.Item(2).RowHeight = 27
$ws.Rows.Item(3).RowHeight = 27

# Row 4's "No" counter is no longer needed now that rows 2-3 share one.
$ws.Range("A4").ClearContents()

# Leave the selection where the author left it after making the edit.
$ws.Range("A4").Select()
